$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 138, pushing the existing row 138 (and all rows below it,
# through the former row 167) down by one. This matches the diff, where the
# former rows 138-167 become rows 139-168, and a brand new data row is
# introduced at row 138.
$ws.Rows.Item(138).Insert()

# Populate the new row 138 with the new weekly price-survey record.
$ws.Range("A138").Value() = 4
$ws.Range("B138").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C138").Value() = "Los Lagos"
$ws.Range("D138").Value() = 44476
$ws.Range("E138").Value() = 10
$ws.Range("F138").Value() = 100112037
$ws.Range("G138").Value() = "Cebollín"
$ws.Range("H138").Value() = "Sin especificar"
$ws.Range("I138").Value() = "Primera"
$ws.Range("J138").Value() = 80
$ws.Range("K138").Value() = 6000
$ws.Range("L138").Value() = 6000
$ws.Range("M138").Value() = 6000
$ws.Range("N138").Value() = "$/paquete 36 unidades"
$ws.Range("O138").Value() = "Región Metropolitana"
$ws.Range("P138").Value() = 167
$ws.Range("Q138").Value() = 36
$ws.Range("R138").Value() = "Hortaliza"

# Ensure the date cell keeps the date/time number format used by the rest of
# the "Fecha" column.
$ws.Range("D138").NumberFormat = $ws.Range("D139").NumberFormat
